$wb = $excel.ActiveWorkbook

# ---- TODOS sheet ----
$wsTodos = $wb.Worksheets.Item("TODOS")
$todosData = @(
  @("16:14", "215C_LA PLATA", 0, "📅"),
  @("16:15", "16_SANTA ANA", 1, "🚌"),
  @("16:15", "16_SANTA ANA", 1, "🚌"),
  @("16:17", "225_C ROCA-H SUR", 3, "🚌"),
  @("16:20", "215C_EL PATO", 6, "📅"),
  @("16:21", "26_HERNANDEZ", 7, "📅"),
  @("16:29", "10_OLMOS", 15, "🚌"),
  @("16:30", "15_ABASTO", 16, "🚌"),
  @("16:36", "11_ETCHEVERRY", 22, "🚌"),
  @("16:40", "17_ROMERO", 26, "📅"),
  @("16:42", "16_P MOR-SANTA ANA", 28, "🚌"),
  @("16:43", "225_GOMEZ", 29, "📅"),
  @("16:48", "15_ABASTO", 34, "🚌"),
  @("16:50", "14_ABASTO", 36, "🚌"),
  @("16:53", "215B_LP-P MOR-40 Y 115", 39, "🚌"),
  @("16:56", "17_179 Y 38", 42, "🚌"),
  @("16:57", "10_OLMOS", 43, "🚌"),
  @("17:04", "11_ETCHEVERRY", 50, "🚌"),
  @("17:04", "215A_EL PATO", 50, "🚌"),
  @("17:04", "23_HERNANDEZ", 50, "🚌"),
  @("17:14", "215A_LA PLATA", 60, "🚌"),
  @("17:21", "26_HERNANDEZ", 67, "🚌"),
  @("17:28", "14_ABASTO", 74, "🚌"),
  @("17:40", "215B_EL PATO", 86, "📅"),
  @("17:41", "84_COLONIA URQUIZA-ESC 49", 87, "🚌"),
  @("17:43", "27_EL RETIRO", 89, "🚌"),
  @("17:44", "23_HERNANDEZ", 90, "🚌"),
  @("17:50", "16_P MOR-167 Y 521", 96, "📅"),
  @("17:52", "81_EL PELIGRO", 98, "🚌"),
  @("18:04", "17_ROMERO", 110, "🚌"),
  @("18:04", "215C_LA PLATA", 110, "📅")
)
for ($i = 0; $i -lt $todosData.Length; $i++) {
  $row = $todosData[$i]
  $r = $i + 2
  $wsTodos.Cells.Item($r, 1).Value = $row[0]
  $wsTodos.Cells.Item($r, 2).Value = $row[1]
  $wsTodos.Cells.Item($r, 3).Value = $row[2]
  $wsTodos.Cells.Item($r, 4).Value = $row[3]
}

# ---- 215 sheet ----
$ws215 = $wb.Worksheets.Item("215")
$s215Data = @(
  @("16:14", "215C_LA PLATA", 0, "📅"),
  @("16:20", "215C_EL PATO", 6, "📅"),
  @("16:53", "215B_LP-P MOR-40 Y 115", 39, "🚌"),
  @("17:04", "215A_EL PATO", 50, "🚌"),
  @("17:14", "215A_LA PLATA", 60, "🚌"),
  @("17:40", "215B_EL PATO", 86, "📅"),
  @("18:04", "215C_LA PLATA", 110, "📅")
)
for ($i = 0; $i -lt $s215Data.Length; $i++) {
  $row = $s215Data[$i]
  $r = $i + 2
  $ws215.Cells.Item($r, 1).Value = $row[0]
  $ws215.Cells.Item($r, 2).Value = $row[1]
  $ws215.Cells.Item($r, 3).Value = $row[2]
  $ws215.Cells.Item($r, 4).Value = $row[3]
}

# ---- COMBINADAS sheet (same data as TODOS) ----
$wsComb = $wb.Worksheets.Item("COMBINADAS")
for ($i = 0; $i -lt $todosData.Length; $i++) {
  $row = $todosData[$i]
  $r = $i + 2
  $wsComb.Cells.Item($r, 1).Value = $row[0]
  $wsComb.Cells.Item($r, 2).Value = $row[1]
  $wsComb.Cells.Item($r, 3).Value = $row[2]
  $wsComb.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "Done. TODOS dim: $($wsTodos.UsedRange.Address()) 215 dim: $($ws215.UsedRange.Address()) COMBINADAS dim: $($wsComb.UsedRange.Address())"
